$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 4 (even_MAG-GUT3315.fa). Remaining rows shift up, so the
# former row 5 (even_MAG-GUT81714.fa) becomes the new row 4 and the
# former row 6 (even_MAG-GUT81784.fa) becomes the new row 5.
$ws.Rows("4").Delete()

# Remove the new row 5 (originally row 6, even_MAG-GUT81784.fa),
# leaving the sheet with data rows 1-4 only.
$ws.Rows("5").Delete()
